$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 - only Taxonsorteringsordning (B) changes
$ws.Range("B30").Value = 89936

# Row 31 - becomes the old row-33 species data (Lakritsmusseron) with new coords
$ws.Range("A31").Value = 112017534
$ws.Range("B31").Value = 88126
$ws.Range("D31").Value = "VU"
$ws.Range("E31").Value = 1593
$ws.Range("F31").Value = "Lakritsmusseron"
$ws.Range("G31").Value = "Tricholoma apium"
$ws.Range("H31").Value = "Jul.Schäff."
$ws.Range("I31").NumberFormat = "@"
$ws.Range("I31").Value = "4"
$ws.Range("I31").NumberFormat = "General"
$ws.Range("J31").Value = "fruktkroppar"
$ws.Range("Q31").Value = 683073
$ws.Range("R31").Value = 6575478

# Row 32 - becomes the old row-34 species data (Dropptaggsvamp) with new coords
$ws.Range("A32").Value = 112017447
$ws.Range("B32").Value = 90800
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 4364
$ws.Range("F32").Value = "Dropptaggsvamp"
$ws.Range("G32").Value = "Hydnellum ferrugineum"
$ws.Range("H32").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q32").Value = 682844
$ws.Range("R32").Value = 6575514

# Row 33 - becomes the old row-31 species data (Svartvit taggsvamp) with new coords
$ws.Range("A33").Value = 112017413
$ws.Range("B33").Value = 90843
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 5448
$ws.Range("F33").Value = "Svartvit taggsvamp"
$ws.Range("G33").Value = "Phellodon connatus"
$ws.Range("H33").Value = "(Schultz) nom.prov"
$ws.Range("I33").Value = ""
$ws.Range("J33").Value = ""
$ws.Range("Q33").Value = 682734
$ws.Range("R33").Value = 6575482

# Row 34 - keeps Dropptaggsvamp species, just Id/order/coords change
$ws.Range("A34").Value = 112017252
$ws.Range("B34").Value = 90800
$ws.Range("Q34").Value = 682711
$ws.Range("R34").Value = 6575494

# Row 35 - becomes the old row-36 species data (Skarp dropptaggsvamp) with new coords
$ws.Range("A35").Value = 112017488
$ws.Range("B35").Value = 90812
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 4366
$ws.Range("F35").Value = "Skarp dropptaggsvamp"
$ws.Range("G35").Value = "Hydnellum peckii"
$ws.Range("H35").Value = "Banker"
$ws.Range("I35").Value = ""
$ws.Range("J35").Value = ""
$ws.Range("Q35").Value = 682956
$ws.Range("R35").Value = 6575474

# Row 36 - becomes the old row-32 species data (Svart taggsvamp) with new coords
$ws.Range("A36").Value = 112017392
$ws.Range("B36").Value = 90844
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 5449
$ws.Range("F36").Value = "Svart taggsvamp"
$ws.Range("G36").Value = "Phellodon niger"
$ws.Range("H36").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q36").Value = 682712
$ws.Range("R36").Value = 6575458

# Row 37 - becomes the old row-35 species data (Goliatmusseron) with new coords
$ws.Range("A37").Value = 112017512
$ws.Range("B37").Value = 88166
$ws.Range("D37").Value = "VU"
$ws.Range("E37").Value = 6276
$ws.Range("F37").Value = "Goliatmusseron"
$ws.Range("G37").Value = "Tricholoma matsutake"
$ws.Range("H37").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "4"
$ws.Range("I37").NumberFormat = "General"
$ws.Range("J37").Value = "fruktkroppar"
$ws.Range("Q37").Value = 683037
$ws.Range("R37").Value = 6575484

# Row 38 - becomes the old row-41/36-style Skarp dropptaggsvamp data with new Id/coords
$ws.Range("A38").Value = 112017224
$ws.Range("B38").Value = 90812
$ws.Range("E38").Value = 4366
$ws.Range("F38").Value = "Skarp dropptaggsvamp"
$ws.Range("G38").Value = "Hydnellum peckii"
$ws.Range("H38").Value = "Banker"
$ws.Range("Q38").Value = 682703
$ws.Range("R38").Value = 6575491

# Row 39 - becomes Svart taggsvamp data with new Id/coords
$ws.Range("A39").Value = 112017159
$ws.Range("B39").Value = 90844
$ws.Range("E39").Value = 5449
$ws.Range("F39").Value = "Svart taggsvamp"
$ws.Range("G39").Value = "Phellodon niger"
$ws.Range("H39").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q39").Value = 682699
$ws.Range("R39").Value = 6575482

# Row 40 - only Taxonsorteringsordning (B) changes
$ws.Range("B40").Value = 88166

# Row 41 - becomes Svartvit taggsvamp data with new Id/coords
$ws.Range("A41").Value = 112017430
$ws.Range("B41").Value = 90843
$ws.Range("D41").Value = "NT"
$ws.Range("E41").Value = 5448
$ws.Range("F41").Value = "Svartvit taggsvamp"
$ws.Range("G41").Value = "Phellodon connatus"
$ws.Range("H41").Value = "(Schultz) nom.prov"
$ws.Range("Q41").Value = 682793
$ws.Range("R41").Value = 6575520

# Row 42 - becomes Dropptaggsvamp data with new Id/coords
$ws.Range("A42").Value = 112017130
$ws.Range("B42").Value = 90800
$ws.Range("D42").Value = "LC"
$ws.Range("E42").Value = 4364
$ws.Range("F42").Value = "Dropptaggsvamp"
$ws.Range("G42").Value = "Hydnellum ferrugineum"
$ws.Range("H42").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q42").Value = 682695
$ws.Range("R42").Value = 6575454

# Row 43 - only Taxonsorteringsordning (B) changes
$ws.Range("B43").Value = 90794
